$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Fill in the yearly coverage values for row 2 across all year columns (H:AD = 2018..2040)
# Previously only every other year (H, J, L, N, P, R, T, V, X, Z, AB, AD) had 0.6;
# now every year column should have 0.6.
$ws.Range("H2:AD2").Value = 0.6
